$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while forcing text interpretation
# (so numeric-looking strings like "1.01" are not coerced to numbers),
# then drop back to the default "Normal" style so no stray number format
# is left behind on the cell.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '43.071.31'
$ws.Range("E2").Value = '  -0.08%  '

Set-TextValue $ws.Range("D3") '2.301.82'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue $ws.Range("D5") '301.02'
$ws.Range("E5").Value = '  +0.13%  '

Set-TextValue $ws.Range("D6") '99.62'
$ws.Range("E6").Value = '  +2.27%  '

Set-TextValue $ws.Range("D7") '0.505'
$ws.Range("E7").Value = '  -0.68%  '

$ws.Range("E8").Value = '  -0.02%  '

Set-TextValue $ws.Range("D9") '0.511'
$ws.Range("E9").Value = '  +2.45%  '

Set-TextValue $ws.Range("D10") '36.33'
$ws.Range("E10").Value = '  +8.05%  '

Set-TextValue $ws.Range("D11") '0.0790'
$ws.Range("E11").Value = '  -0.87%  '

Set-TextValue $ws.Range("D12") '0.117'
$ws.Range("E12").Value = '  +1.14%  '

$ws.Range("E13").Value = '  +7.52%  '

$ws.Range("E14").Value = '  +2.28%  '

Set-TextValue $ws.Range("D15") '2.655.90'
$ws.Range("E15").Value = '  -0.27%  '

Set-TextValue $ws.Range("D16") '2.261.91'
$ws.Range("E16").Value = '  -2.96%  '

Set-TextValue $ws.Range("D17") '0.800'
$ws.Range("E17").Value = '  -1.11%  '

Set-TextValue $ws.Range("D18") '42.973.70'

Set-TextValue $ws.Range("D19") '12.59'
$ws.Range("E19").Value = '  +9.08%  '

Set-TextValue $ws.Range("D20") '0.0₃0904'
$ws.Range("E20").Value = '  +0.34%  '

Set-TextValue $ws.Range("D21") '6.13'
$ws.Range("E21").Value = '  +1.35%  '

Set-TextValue $ws.Range("D22") '67.89'
$ws.Range("E22").Value = '  +0.51%  '

Set-TextValue $ws.Range("D23") '236.03'
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("E24").Value = '  +10.56%  '

Set-TextValue $ws.Range("D25") '1.01'
$ws.Range("E25").Value = '  +0.69%  '

Set-TextValue $ws.Range("D26") '2.45'
$ws.Range("E26").Value = '  -0.41%  '

Set-TextValue $ws.Range("D27") '24.99'
$ws.Range("E27").Value = '  +2.34%  '

Set-TextValue $ws.Range("D28") '2.36'
$ws.Range("E28").Value = '  +9.19%  '

Set-TextValue $ws.Range("D29") '34.65'
$ws.Range("E29").Value = '  +2.11%  '

Set-TextValue $ws.Range("D30") '167.16'
$ws.Range("E30").Value = '  +0.25%  '

Set-TextValue $ws.Range("D31") '9.14'
$ws.Range("E31").Value = '  -0.03%  '

$ws.Range("E32").Value = '  -0.05%  '

Set-TextValue $ws.Range("D33") '5.03'
$ws.Range("E33").Value = '  +1.61%  '

Set-TextValue $ws.Range("D34") '17.72'
$ws.Range("E34").Value = '  +4.56%  '

Set-TextValue $ws.Range("D35") '4.61'
$ws.Range("E35").Value = '  -0.72%  '

$ws.Range("E36").Value = '  +1.10%  '

$ws.Range("E37").Value = '  -0.37%  '

# Row 38: coin ranking swapped in with LidoDAOToken
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D38") '2.82'
$ws.Range("E38").Value = '  -0.04%  '

# Row 39: coin ranking swapped in with ARBITRUM
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range("D39") '1.79'
$ws.Range("E39").Value = '  +2.05%  '

$ws.Range("E40").Value = '  -0.62%  '

Set-TextValue $ws.Range("D41") '0.110'
$ws.Range("E41").Value = '  +0.06%  '

Set-TextValue $ws.Range("D42") '2.31'
$ws.Range("E42").Value = '  -1.18%  '

$ws.Range("E43").Value = '  +3.40%  '

Set-TextValue $ws.Range("D44") '1.978.51'
$ws.Range("E44").Value = '  -0.35%  '

$ws.Range("E45").Value = '  +3.08%  '

Set-TextValue $ws.Range("D46") '2.91'
$ws.Range("E46").Value = '  +1.87%  '

Set-TextValue $ws.Range("D47") '17.54'
$ws.Range("E47").Value = '  -0.21%  '

Set-TextValue $ws.Range("D48") '55.64'
$ws.Range("E48").Value = '  +4.66%  '

$ws.Range("E49").Value = '  +3.76%  '

Set-TextValue $ws.Range("D50") '2.522.81'
$ws.Range("E50").Value = '  -0.22%  '

Set-TextValue $ws.Range("D51") '70.84'
$ws.Range("E51").Value = '  +1.10%  '
